$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7142412066459656
$ws.Range("B1").Value = 3.587246179580688
$ws.Range("C1").Value = 4.299013137817383
$ws.Range("D1").Value = 2.624518871307373
$ws.Range("E1").Value = 1.058185338973999
